$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.031.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.81%  '

$ws.Range("D3").Value = '''1.808.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.32%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.73%  '

$ws.Range("D5").Value = '''329.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '

$ws.Range("D7").Value = '''0.4434'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.77%  '

$ws.Range("D8").Value = '''0.3724'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.82%  '

$ws.Range("D9").Value = '''44.67'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.51%  '

$ws.Range("D10").Value = '''0.07697'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.24%  '

$ws.Range("E11").Value = '  -0.90%  '

$ws.Range("D12").Value = '''1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.68%  '

$ws.Range("D13").Value = '''21.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").Value = '''6.290'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").Value = '''7.460'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").Value = '''1.816.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.74%  '

$ws.Range("D17").Value = '''93.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +12.07%  '

$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("E19").Value = '  -1.22%  '

$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").Value = '''17.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '

$ws.Range("D22").Value = '''6.255'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").Value = '''0.5347'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.06%  '

$ws.Range("D24").Value = '''28.072.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.76%  '

$ws.Range("E25").Value = '  +2.92%  '

$ws.Range("D26").Value = '''2.114'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -12.13%  '

$ws.Range("D27").Value = '''20.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.51%  '

$ws.Range("E28").Value = '  +2.14%  '

$ws.Range("D29").Value = '''2.020.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.41%  '

$ws.Range("D30").Value = '''2.320'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.06%  '

$ws.Range("D31").Value = '''126.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("D32").Value = '''1.202'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.23%  '

$ws.Range("D33").Value = '''5.846'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.21%  '

$ws.Range("D34").Value = '''0.09232'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.49%  '

$ws.Range("D35").Value = '''3.666'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.84%  '

$ws.Range("D36").Value = '''13.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.35%  '

$ws.Range("D37").Value = '''0.02344'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.40%  '

$ws.Range("D38").Value = '''0.2168'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.46%  '

$ws.Range("D39").Value = '''5.161'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6561'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '''0.06199'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.29%  '

$ws.Range("D42").Value = '''1.194'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = '''8.087'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("D44").Value = '''1.003'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.63%  '

$ws.Range("D45").Value = '''13.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.48%  '

$ws.Range("E46").Value = '  -3.30%  '

$ws.Range("E47").Value = '  +1.47%  '

$ws.Range("D48").Value = '''3.763'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.78%  '

$ws.Range("D49").Value = '''126.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").Value = '''2.029'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.73%  '

$ws.Range("D51").Value = '''1.151'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.85%  '
